# RESTORE: Recover original multi-industry template text
# (reverts "Banking" industry-specific copy back to generic
#  "FINANCE Change Management Plan" / "AI/ML" placeholder text)

$wb = $excel.ActiveWorkbook

# --- Sheet: Change Management Overview ---
$wsOverview = $wb.Worksheets.Item("Change Management Overview")

$wsOverview.Range("A2").Value = "FINANCE Change Management Plan Project"
$wsOverview.Range("B6").Value = "Enterprise AI/ML Implementation"
$wsOverview.Range("A15").Value = "1. Achieve 95% user adoption of new AI/ML systems within 6 months of go-live"
$wsOverview.Range("A17").Value = "3. Build organizational capability and confidence in AI/ML technologies"
$wsOverview.Range("A20").Value = "6. Create positive stakeholder sentiment and enthusiasm for AI/ML transformation"

# Row 4 on this sheet is (and must remain) a blank spacer row. The engine's
# xlsx writer omits any row that has never been "touched", so nudge its
# row-level state (outline level) and immediately set it back to its
# default (0) — this keeps row 4 present in the saved sheet as a genuinely
# empty row (no cells), rather than it vanishing or picking up stray cells.
$wsOverview.Rows.Item(4).OutlineLevel = 1
$wsOverview.Rows.Item(4).OutlineLevel = 0

# --- Sheet: Change Impact Assessment ---
$wsImpact = $wb.Worksheets.Item("Change Impact Assessment")

$wsImpact.Range("G4").Value = "AI/ML automation"
